$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 3

# Delete row 3 entirely (shifts cells up, removes the row)
$ws.Range("A3:B3").Delete()
